$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume/coin data per latest scrape
$ws.Range("D2").Value = '58.842.94'
$ws.Range("E2").Value = '  -3.51%  '
$ws.Range("D3").Value = '3.207.41'
$ws.Range("E3").Value = '  -4.78%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.35%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '3.207.67'
$ws.Range("E8").Value = '  -4.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.457'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.51'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.82%  '
$ws.Range("E11").Value = '  -6.82%  '
$ws.Range("E12").Value = '  -5.82%  '
$ws.Range("D13").Value = '3.759.73'
$ws.Range("E13").Value = '  -4.82%  '
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.98%  '
$ws.Range("D16").Value = '3.210.96'
$ws.Range("E16").Value = '  -4.64%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '58.873.74'
$ws.Range("E17").Value = '  -3.58%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000157'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.43%  '
$ws.Range("E19").Value = '  -7.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '360.58'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.95%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.14%  '
$ws.Range("E25").Value = '  -7.72%  '
$ws.Range("D26").Value = '3.343.00'
$ws.Range("E26").Value = '  -4.64%  '
$ws.Range("E27").Value = '  -2.99%  '
$ws.Range("D28").Value = '0.0₃0963'
$ws.Range("E28").Value = '  -11.18%  '
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.74'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.68%  '
$ws.Range("E37").Value = '  -9.04%  '
$ws.Range("E38").Value = '  -6.76%  '
$ws.Range("E39").Value = '  -8.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0705'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.59%  '
$ws.Range("D42").Value = '3.236.37'
$ws.Range("E42").Value = '  -4.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.713'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.42%  '
$ws.Range("E45").Value = '  -4.19%  '
$ws.Range("E46").Value = '  -6.70%  '
$ws.Range("E47").Value = '  -7.08%  '
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").Value = '2.302.25'
$ws.Range("E49").Value = '  -7.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.56%  '
